# Extent Report Screenshot Update Issue Fix
#
# The TestSuite's RunMode flags for the "AddCustomerTest" and
# "OpenAccountTest" rows were set to "N", which skipped those tests (and
# therefore the Extent Report screenshot capture for them). Flip both
# RunMode cells to "Y" so every test actually runs, and move the active
# sheet/selection back to the TestSuite tab.

$wb = $excel.ActiveWorkbook

$testSuite = $wb.Worksheets.Item("TestSuite")
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")

# TestSuite!B4 (OpenAccountTest's RunMode): N -> Y
$testSuite.Range("B4").Value = "Y"

# AddCustomerTest!E3 (row 3's RunMode): N -> Y
$addCustomer.Range("E3").Value = "Y"

# Normalize the formatting of the RunMode column on AddCustomerTest so it
# matches the rest of the sheet (copy the plain bordered style from D2,
# formats only, onto E2:E4 - values are left untouched).
$addCustomer.Range("D2").Copy()
$addCustomer.Range("E2:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection on AddCustomerTest to D7 for the next time it's
# viewed, then switch the active tab back to TestSuite with B4 selected.
[void]$addCustomer.Range("D7").Select()

[void]$testSuite.Activate()
[void]$testSuite.Range("B4").Select()
